# The document contains three occurrences of an "<id>...</id>" marker,
# each split across three runs:
#   run 1: "<id>"   (Courier New, color 7f6000, sz 18)
#   run 2: "p124v_N" (color 000000, plain)
#   run 3: "</id>"  (Courier New, color 7f6000, sz 18 - same formatting as run 1)
#
# The edit merges each triple into a single run holding the full literal
# text "<id>p124v_N</id>", using run 1's (Courier New / 7f6000) formatting,
# and removes the now-redundant runs 2 and 3. Word's Find/Replace naturally
# collapses the matched runs into one run carrying the formatting of the
# first character of the match, so a plain literal Find.Execute achieves
# exactly that merge for each of the three ids (p124v_1, p124v_2, p124v_3).
# The unrelated "fig_p124v_N" markers elsewhere in the document are left
# untouched since the search text is an exact, literal match.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p124v_1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p124v_1</id>", 2)
$d.Content.Find.Execute("<id>p124v_2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p124v_2</id>", 2)
$d.Content.Find.Execute("<id>p124v_3</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p124v_3</id>", 2)
